$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin price list refresh (cryptos.xlsx) -- values/percentages updated per latest pull.
# Column D ("Price") cells are forced back to Text after the write so purely-numeric-
# looking prices (e.g. "567.64") do not get auto-promoted to the Number type by Excel,
# matching the inline-string cells already used throughout this sheet.

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.466.42'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.49%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.541.88'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +4.74%  '

# Row 4
$ws.Range('E4').Value = '  +0.07%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '567.64'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.88%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '150.36'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +8.24%  '

# Row 7
$ws.Range('E7').Value = '  +0.07%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.585'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.02%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.540.13'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +4.73%  '

# Row 10
$ws.Range('E10').Value = '  +1.58%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.69'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.56%  '

# Row 12
$ws.Range('E12').Value = '  +1.03%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.357'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.63%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.18'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +7.71%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.997.60'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.81%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.397.32'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.62%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000143'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.74%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.539.56'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +4.60%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.57'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.10%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '338.93'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.95%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.33'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.83%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.78'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.92%  '

# Row 23
$ws.Range('E23').Value = '  +0.10%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '65.93'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.01%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.169'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.73%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.55'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +14.67%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.59'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.66%  '

# Row 28
$ws.Range('B28').Value = 'Binance-PegBSC-USD'
$ws.Range('C28').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.29%  '

# Row 29
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.44'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.49%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.14'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +10.64%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0₃0816'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.13%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.86'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.40%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '177.75'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.17%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.58'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +9.49%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '419.38'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +10.83%  '

# Row 36
$ws.Range('E36').Value = '  +2.36%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '19.01'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.42%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.42'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.70%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.76'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.75%  '

# Row 41
$ws.Range('E41').Value = '  -0.03%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '39.49'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.38%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '153.44'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +5.88%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.77'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.76%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '20.71'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.42%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.608'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.63%  '

# Row 47
$ws.Range('B47').Value = 'Stellar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0964'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.17%  '

# Row 48
$ws.Range('B48').Value = 'Hedera'
$ws.Range('C48').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0524'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.81%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0238'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +6.68%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '18.63'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.38%  '

# Row 51
$ws.Range('E51').Value = '  +5.11%  '
